# "Updated test case results"
#
# The Multiplicative sheet's comparison columns (K = "VSCode Multi",
# L = "VSCode Match" helper math) were refreshed with new computed
# results. Updating K/L causes the dependent shared formulas in S
# ("VSCode Match") and T to recalculate automatically.
#
# The author was also looking at / working on the "Multiplicative" tab
# when the file was saved (it becomes the active tab, with a fresh
# selection at J26), whereas "Reductive" was the active tab before.

$wb = $excel.ActiveWorkbook
$ws2 = $wb.Worksheets.Item("Multiplicative")

# New K/L results for the affected rows.
$results = @{
    10 = @(5.2031851927965702, 2.8048366862753902)
    11 = @(6.5989391226003402, 4.1339271252336296)
    12 = @(5.2031851927965702, 2.8048366862753902)
    13 = @(5.22945719061173,   2.8197091778474701)
    16 = @(5.2285133096524401, 2.8192002397862699)
    17 = @(6.5975399762486298, 4.1330506253981198)
    18 = @(5.2285133096524401, 2.8192002397862699)
    21 = @(5.2294571793469604, 2.8197091717735399)
    22 = @(6.5989391059021898, 4.1339271147730203)
    23 = @(5.2294571793469604, 2.8197091717735399)
    29 = @(7.1226756462811602, 3.9115964794254401)
    30 = @(10.934154886446199, 6.9764918155464004)
    31 = @(7.1226756462811602, 3.9115964794254401)
    35 = @(7.1226756639777102, 3.9115964891439399)
    36 = @(10.934154917711499, 6.9764918354950902)
    37 = @(7.1225465300232997, 3.9115178901755598)
    42 = @(10.934154886446199, 6.9764918155464004)
    43 = @(7.1226756462811602, 3.9115964794254401)
    49 = @(7.1030022516433702, 3.89962895042376)
}

foreach ($row in $results.Keys) {
    $vals = $results[$row]
    $ws2.Range("K$row").Value = $vals[0]
    $ws2.Range("L$row").Value = $vals[1]
}

# Switch focus to the Multiplicative sheet, matching the selection left
# behind in the saved file.
$ws2.Activate() | Out-Null
$ws2.Range("J26").Select() | Out-Null
